# Auto-generated: apply numeric value updates per the authoritative diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4630.4375
$ws.Cells.Item(32, 9).Value = 11367.75
$ws.Cells.Item(32, 10).Value = 2384.6667
$ws.Cells.Item(32, 11).Value = 11367.75
$ws.Cells.Item(32, 12).Value = 2384.6667
$ws.Cells.Item(32, 13).Value = -11041.75
$ws.Cells.Item(32, 14).Value = -3036.6667
$ws.Cells.Item(40, 8).Value = 8813
$ws.Cells.Item(40, 10).Value = 13126
$ws.Cells.Item(40, 12).Value = 13126
$ws.Cells.Item(40, 14).Value = -13476
$ws.Cells.Item(41, 8).Value = 693.625
$ws.Cells.Item(41, 10).Value = 1087.25
$ws.Cells.Item(41, 12).Value = 1087.25
$ws.Cells.Item(41, 14).Value = -1967.25
$ws.Cells.Item(98, 8).Value = 316647.7
$ws.Cells.Item(98, 10).Value = 1000351
$ws.Cells.Item(98, 12).Value = 1000351
$ws.Cells.Item(98, 14).Value = -1003347
$ws.Cells.Item(122, 8).Value = 316647.7
$ws.Cells.Item(122, 10).Value = 1000351
$ws.Cells.Item(122, 12).Value = 3001053
$ws.Cells.Item(122, 14).Value = -3005953
$ws.Cells.Item(136, 8).Value = 69999
$ws.Cells.Item(136, 10).Value = 69999
$ws.Cells.Item(136, 12).Value = 69999
$ws.Cells.Item(136, 14).Value = -80199

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8856.923000000001
$ws.Cells.Item(32, 9).Value = 7512.7
$ws.Cells.Item(32, 11).Value = 7512.7
$ws.Cells.Item(32, 13).Value = -7225.7
$ws.Cells.Item(63, 8).Value = 6455.8887
$ws.Cells.Item(63, 9).Value = 2025.75
$ws.Cells.Item(63, 10).Value = 10000
$ws.Cells.Item(63, 11).Value = 2025.75
$ws.Cells.Item(63, 12).Value = 10000
$ws.Cells.Item(63, 13).Value = -1339.75
$ws.Cells.Item(63, 14).Value = -11372
$ws.Cells.Item(66, 8).Value = 6455.8887
$ws.Cells.Item(66, 9).Value = 2025.75
$ws.Cells.Item(66, 10).Value = 10000
$ws.Cells.Item(66, 11).Value = 10128.75
$ws.Cells.Item(66, 12).Value = 50000
$ws.Cells.Item(66, 13).Value = -6696.75
$ws.Cells.Item(66, 14).Value = -56864
$ws.Cells.Item(112, 8).Value = 39535.57
$ws.Cells.Item(112, 10).Value = 39535.57
$ws.Cells.Item(112, 12).Value = 39535.57
$ws.Cells.Item(112, 14).Value = -42489.57
$ws.Cells.Item(132, 8).Value = 3148
$ws.Cells.Item(132, 9).Value = 2090.2666
$ws.Cells.Item(132, 11).Value = 6270.7998
$ws.Cells.Item(132, 13).Value = -3740.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(102, 8).Value = 4400
$ws.Cells.Item(102, 9).Value = 4400
$ws.Cells.Item(102, 11).Value = 4400
$ws.Cells.Item(102, 13).Value = -1155
$ws.Cells.Item(140, 8).Value = 51611.555
$ws.Cells.Item(140, 10).Value = 51611.555
$ws.Cells.Item(140, 12).Value = 51611.555
$ws.Cells.Item(140, 14).Value = -61971.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 1433.6666
$ws.Cells.Item(6, 9).Value = 1499.5
$ws.Cells.Item(6, 10).Value = 1302
$ws.Cells.Item(6, 11).Value = 1499.5
$ws.Cells.Item(6, 12).Value = 1302
$ws.Cells.Item(6, 13).Value = -1386.5
$ws.Cells.Item(6, 14).Value = -1528
$ws.Cells.Item(16, 8).Value = 2543.4375
$ws.Cells.Item(16, 9).Value = 2131.6
$ws.Cells.Item(16, 10).Value = 3229.8333
$ws.Cells.Item(16, 11).Value = 2131.6
$ws.Cells.Item(16, 12).Value = 3229.8333
$ws.Cells.Item(16, 13).Value = -1844.6
$ws.Cells.Item(16, 14).Value = -3803.8333
$ws.Cells.Item(31, 8).Value = 33972.73
$ws.Cells.Item(31, 9).Value = 3816.158
$ws.Cells.Item(31, 10).Value = 65804.664
$ws.Cells.Item(31, 11).Value = 3816.158
$ws.Cells.Item(31, 12).Value = 65804.664
$ws.Cells.Item(31, 13).Value = -3521.158
$ws.Cells.Item(31, 14).Value = -66394.664
$ws.Cells.Item(34, 8).Value = 33972.73
$ws.Cells.Item(34, 9).Value = 3816.158
$ws.Cells.Item(34, 10).Value = 65804.664
$ws.Cells.Item(34, 11).Value = 3816.158
$ws.Cells.Item(34, 12).Value = 65804.664
$ws.Cells.Item(34, 13).Value = -3614.158
$ws.Cells.Item(34, 14).Value = -66208.664
$ws.Cells.Item(99, 8).Value = 3299.4285
$ws.Cells.Item(99, 9).Value = 3049.2
$ws.Cells.Item(99, 10).Value = 3925
$ws.Cells.Item(99, 11).Value = 3049.2
$ws.Cells.Item(99, 12).Value = 3925
$ws.Cells.Item(99, 13).Value = -1551.2
$ws.Cells.Item(99, 14).Value = -6921
$ws.Cells.Item(113, 8).Value = 2543.4375
$ws.Cells.Item(113, 9).Value = 2131.6
$ws.Cells.Item(113, 10).Value = 3229.8333
$ws.Cells.Item(113, 11).Value = 2131.6
$ws.Cells.Item(113, 12).Value = 3229.8333
$ws.Cells.Item(113, 13).Value = 38.40000000000009
$ws.Cells.Item(113, 14).Value = -7569.8333
$ws.Cells.Item(122, 8).Value = 8409.5
$ws.Cells.Item(122, 9).Value = 5283.6665
$ws.Cells.Item(122, 11).Value = 15850.9995
$ws.Cells.Item(122, 13).Value = -13400.9995
$ws.Cells.Item(126, 8).Value = 3299.4285
$ws.Cells.Item(126, 9).Value = 3049.2
$ws.Cells.Item(126, 10).Value = 3925
$ws.Cells.Item(126, 11).Value = 9147.599999999999
$ws.Cells.Item(126, 12).Value = 11775
$ws.Cells.Item(126, 13).Value = -6677.599999999999
$ws.Cells.Item(126, 14).Value = -16715
$ws.Cells.Item(134, 8).Value = 2639.611
$ws.Cells.Item(134, 9).Value = 2088.0908
$ws.Cells.Item(134, 11).Value = 6264.2724
$ws.Cells.Item(134, 13).Value = -3729.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 1450.9286
$ws.Cells.Item(3, 9).Value = 1186
$ws.Cells.Item(3, 11).Value = 3558
$ws.Cells.Item(3, 13).Value = -3446
$ws.Cells.Item(14, 8).Value = 2885.3333
$ws.Cells.Item(14, 9).Value = 2885.3333
$ws.Cells.Item(14, 11).Value = 8655.999899999999
$ws.Cells.Item(14, 13).Value = -8482.999899999999
$ws.Cells.Item(23, 8).Value = 190.66667
$ws.Cells.Item(23, 9).Value = 123.5
$ws.Cells.Item(23, 10).Value = 209.85715
$ws.Cells.Item(23, 11).Value = 370.5
$ws.Cells.Item(23, 12).Value = 629.5714499999999
$ws.Cells.Item(23, 13).Value = -135.5
$ws.Cells.Item(23, 14).Value = -1099.57145
$ws.Cells.Item(68, 8).Value = 20834802
$ws.Cells.Item(68, 9).Value = 62500310
$ws.Cells.Item(68, 11).Value = 187500930
$ws.Cells.Item(68, 13).Value = -187500119
$ws.Cells.Item(71, 8).Value = 20834802
$ws.Cells.Item(71, 9).Value = 62500310
$ws.Cells.Item(71, 11).Value = 562502790
$ws.Cells.Item(71, 13).Value = -562498734
$ws.Cells.Item(98, 8).Value = 1342.6111
$ws.Cells.Item(98, 9).Value = 1002.55554
$ws.Cells.Item(98, 10).Value = 1682.6666
$ws.Cells.Item(98, 11).Value = 3007.66662
$ws.Cells.Item(98, 12).Value = 5047.9998
$ws.Cells.Item(98, 13).Value = -1509.66662
$ws.Cells.Item(98, 14).Value = -8043.9998
$ws.Cells.Item(133, 8).Value = 83335080
$ws.Cells.Item(133, 9).Value = 2326.6667
$ws.Cells.Item(133, 10).Value = 333333340
$ws.Cells.Item(133, 11).Value = 6980.000100000001
$ws.Cells.Item(133, 12).Value = 1000000020
$ws.Cells.Item(133, 13).Value = -1920.000100000001
$ws.Cells.Item(133, 14).Value = -1000010140
$ws.Cells.Item(136, 8).Value = 25644040
$ws.Cells.Item(136, 9).Value = 30305684
$ws.Cells.Item(136, 11).Value = 90917052
$ws.Cells.Item(136, 13).Value = -90911952
$ws.Cells.Item(137, 8).Value = 252050.25
$ws.Cells.Item(137, 9).Value = 2600
$ws.Cells.Item(137, 10).Value = 501500.5
$ws.Cells.Item(137, 11).Value = 7800
$ws.Cells.Item(137, 12).Value = 1504501.5
$ws.Cells.Item(137, 13).Value = -2700
$ws.Cells.Item(137, 14).Value = -1514701.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7349.1577
$ws.Cells.Item(122, 9).Value = 6671.846
$ws.Cells.Item(122, 11).Value = 20015.538
$ws.Cells.Item(122, 13).Value = -17565.538
$ws.Cells.Item(132, 8).Value = 7644.7407
$ws.Cells.Item(132, 9).Value = 7207.4614
$ws.Cells.Item(132, 11).Value = 21622.3842
$ws.Cells.Item(132, 13).Value = -19092.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4929.1904
$ws.Cells.Item(22, 9).Value = 2585
$ws.Cells.Item(22, 10).Value = 9617.571
$ws.Cells.Item(22, 11).Value = 2585
$ws.Cells.Item(22, 12).Value = 9617.571
$ws.Cells.Item(22, 13).Value = -2290
$ws.Cells.Item(22, 14).Value = -10207.571
$ws.Cells.Item(27, 8).Value = 4929.1904
$ws.Cells.Item(27, 9).Value = 2585
$ws.Cells.Item(27, 10).Value = 9617.571
$ws.Cells.Item(27, 11).Value = 2585
$ws.Cells.Item(27, 12).Value = 9617.571
$ws.Cells.Item(27, 13).Value = -2478
$ws.Cells.Item(27, 14).Value = -9831.571
$ws.Cells.Item(54, 8).Value = 64995
$ws.Cells.Item(54, 10).Value = 64995
$ws.Cells.Item(54, 12).Value = 64995
$ws.Cells.Item(54, 14).Value = -66283

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1176.2858
$ws.Cells.Item(100, 9).Value = 848.3125
$ws.Cells.Item(100, 10).Value = 2225.8
$ws.Cells.Item(100, 11).Value = 1696.625
$ws.Cells.Item(100, 12).Value = 4451.6
$ws.Cells.Item(100, 13).Value = -1155.625
$ws.Cells.Item(100, 14).Value = -5533.6
$ws.Cells.Item(136, 8).Value = 3185.5454
$ws.Cells.Item(136, 9).Value = 2432.238
$ws.Cells.Item(136, 11).Value = 7296.714
$ws.Cells.Item(136, 13).Value = -4746.714
